# Add two new columns (I0 / IF) to the sheet, mirroring column H's layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy H1's formatting (bold, centered, bordered) onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values (rows 2-68); column I equals column J for every row.
$values = @(9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,8,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,6,5,6,4,4)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $v = $values[$i]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
